$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows 2-11, columns A-G
# A: Colaborador_id, B: Colaborador_nome, C: Departamento,
# D: Motivo_da_ausência, E: Horas_de_ausência, F: Data_da_ausência, G: Salário

$data = @(
    @{ Row=2;  A=49956; B="Breno Azevedo";        C="TI";                      D="Outros";              E=3; F=45099; G=11023.82 },
    @{ Row=3;  A=49724; B="Emanuella Moreira";     C="Recursos Humanos";        D="Viagem de negócios";  E=1; F=45087; G=7197.09  },
    @{ Row=4;  A=99064; B="Pedro da Mata";         C="Marketing";               D="Outros";              E=8; F=45101; G=5944.47  },
    @{ Row=5;  A=77436; B="Igor Cardoso";          C="Atendimento ao Cliente";  D="Outros";              E=6; F=45090; G=10510.87 },
    @{ Row=6;  A=56215; B="Ana Beatriz Cardoso";   C="Marketing";               D="Outros";              E=4; F=45081; G=3874.15  },
    @{ Row=7;  A=45319; B="Carlos Eduardo Dias";   C="Marketing";               D="Viagem de negócios";  E=1; F=45080; G=6340.72  },
    @{ Row=8;  A=72480; B="Catarina Ramos";        C="Jurídico";                D="Outros";              E=5; F=45092; G=9993.08  },
    @{ Row=9;  A=97058; B="Kaique Rezende";        C="Recursos Humanos";        D="Viagem de negócios";  E=7; F=45100; G=7147.18  },
    @{ Row=10; A=91243; B="Alícia Moraes";         C="P&D";                     D="Doença";              E=1; F=45105; G=11422.72 },
    @{ Row=11; A=33936; B="Yasmin Rocha";          C="Engenharia";              D="Problemas pessoais";  E=5; F=45095; G=2733.06  }
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 1).Value = $item.A
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 3).Value = $item.C
    $ws.Cells.Item($r, 4).Value = $item.D
    $ws.Cells.Item($r, 5).Value = $item.E
    $ws.Cells.Item($r, 6).Value = $item.F
    $ws.Cells.Item($r, 7).Value = $item.G
}

$wb.Save()
